# Replace the field { m:self.name } in the second paragraph with four
# plain-text runs spelling out "{m:self.name" / "" / "" / "}" (the field
# result run keeps its empty middle runs, only the fldChar/instrText
# runs become literal w:t runs), matching TokenIteratorFieldRewriterSplit.

$d = $word.ActiveDocument

$f = $d.Fields.Item(1)

# Remember where the field starts so we can re-insert text there once the
# field (and all of its runs) has been removed.
$fieldStart = $f.Code.Start - 1

# Delete the whole field (fldChar begin/instrText/separate/.../end runs).
# This leaves a clean, empty paragraph behind (its <w:pPr> is untouched).
$f.Delete()

# Build the four-run replacement as inline WordprocessingML and insert it
# at a collapsed range sitting at the (now empty) paragraph's start. Using
# a collapsed range here - rather than the whole paragraph range - makes
# Word merge the new runs into the existing paragraph instead of minting a
# brand new <w:p>, so the original paragraph properties are preserved.
$insertionPoint = $d.Range($fieldStart, $fieldStart)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr>' +
       '<w:r><w:rPr/><w:t>{m:self.name</w:t></w:r>' +
       '<w:r><w:rPr/><w:t/></w:r>' +
       '<w:r><w:rPr/><w:t/></w:r>' +
       '<w:r><w:rPr/><w:t>}</w:t></w:r>' +
       '</w:p>'

$insertionPoint.InsertXML($xml)
